$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 1.966253333333333
$ws.Range("H2").Value = 5.89876
$ws.Range("I2").Value = 0.004409978591445245
$ws.Range("J2").Value = 0.004409978591445245
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 1.105124
$ws.Range("N2").Value = 3.315372
$ws.Range("O2").Value = 0.006910839970832482
$ws.Range("P2").Value = 0.006910839970832482
$ws.Range("Q2").Value = 2.172953748746667
$ws.Range("R2").Value = 19.55658373872
$ws.Range("S2").Value = 0.00003047665632027533
$ws.Range("T2").Value = 0.00003047665632027533

$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 1.966253333333333
$ws.Range("H3").Value = 5.89876
$ws.Range("I3").Value = 0.004409978591445245
$ws.Range("J3").Value = 0.004409978591445245
$ws.Range("O3").Value = 0.0002777950170396876
$ws.Range("P3").Value = 0.0002777950170396876
$ws.Range("Q3").Value = 0.08734621640888889
$ws.Range("R3").Value = 0.78611594768
$ws.Range("S3").Value = 0.000001225070077955189
$ws.Range("T3").Value = 0.000001225070077955189

$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 1.966253333333333
$ws.Range("H4").Value = 5.89876
$ws.Range("I4").Value = 0.004409978591445245
$ws.Range("J4").Value = 0.004409978591445245
$ws.Range("M4").Value = 56.54517366666666
$ws.Range("N4").Value = 169.635521
$ws.Range("O4").Value = 0.3536025335919447
$ws.Range("P4").Value = 0.3536025335919447
$ws.Range("Q4").Value = 111.1821362059955
$ws.Range("R4").Value = 1000.63922585396
$ws.Range("S4").Value = 0.001559379603021274
$ws.Range("T4").Value = 0.001559379603021274

$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 1.966253333333333
$ws.Range("H5").Value = 5.89876
$ws.Range("I5").Value = 0.004409978591445245
$ws.Range("J5").Value = 0.004409978591445245
$ws.Range("M5").Value = 0.8044289999999998
$ws.Range("N5").Value = 2.413287
$ws.Range("O5").Value = 0.005030458199167516
$ws.Range("P5").Value = 0.005030458199167516
$ws.Range("Q5").Value = 1.58171120268
$ws.Range("R5").Value = 14.23540082412
$ws.Range("S5").Value = 0.00002218421296348894
$ws.Range("T5").Value = 0.00002218421296348894

$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 1.966253333333333
$ws.Range("H6").Value = 5.89876
$ws.Range("I6").Value = 0.004409978591445245
$ws.Range("J6").Value = 0.004409978591445245
$ws.Range("M6").Value = 78.08909333333334
$ws.Range("N6").Value = 234.26728
$ws.Range("O6").Value = 0.4883264027331488
$ws.Range("P6").Value = 0.4883264027331488
$ws.Range("Q6").Value = 153.5429400636445
$ws.Range("R6").Value = 1381.8864605728
$ws.Range("S6").Value = 0.002153508981690655
$ws.Range("T6").Value = 0.002153508981690655

$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 1.966253333333333
$ws.Range("H7").Value = 5.89876
$ws.Range("I7").Value = 0.004409978591445245
$ws.Range("J7").Value = 0.004409978591445245
$ws.Range("M7").Value = 23.323433
$ws.Range("N7").Value = 69.970299
$ws.Range("O7").Value = 0.1458519704878668
$ws.Range("P7").Value = 0.1458519704878668
$ws.Range("Q7").Value = 45.85977788102667
$ws.Range("R7").Value = 412.73800092924
$ws.Range("S7").Value = 0.0006432040673715963
$ws.Range("T7").Value = 0.0006432040673715963

$ws.Range("I8").Value = 0.0002834193769566641
$ws.Range("J8").Value = 0.0002834193769566642
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 1.105124
$ws.Range("N8").Value = 3.315372
$ws.Range("O8").Value = 0.006910839970832482
$ws.Range("P8").Value = 0.006910839970832482
$ws.Range("Q8").Value = 0.1396508361333333
$ws.Range("R8").Value = 1.2568575252
$ws.Range("S8").Value = 0.000001958665958780553
$ws.Range("T8").Value = 0.000001958665958780553

$ws.Range("I9").Value = 0.0002834193769566641
$ws.Range("J9").Value = 0.0002834193769566642
$ws.Range("O9").Value = 0.0002777950170396876
$ws.Range("P9").Value = 0.0002777950170396876
$ws.Range("S9").Value = 0.00000007873249065105415
$ws.Range("T9").Value = 0.00000007873249065105416

$ws.Range("I10").Value = 0.0002834193769566641
$ws.Range("J10").Value = 0.0002834193769566642
$ws.Range("M10").Value = 56.54517366666666
$ws.Range("N10").Value = 169.635521
$ws.Range("O10").Value = 0.3536025335919447
$ws.Range("P10").Value = 0.3536025335919447
$ws.Range("Q10").Value = 7.145425112344443
$ws.Range("R10").Value = 64.3088260111
$ws.Range("S10").Value = 0.0001002178097609268
$ws.Range("T10").Value = 0.0001002178097609269

$ws.Range("I11").Value = 0.0002834193769566641
$ws.Range("J11").Value = 0.0002834193769566642
$ws.Range("M11").Value = 0.8044289999999998
$ws.Range("N11").Value = 2.413287
$ws.Range("O11").Value = 0.005030458199167516
$ws.Range("P11").Value = 0.005030458199167516
$ws.Range("Q11").Value = 0.1016530113
$ws.Range("R11").Value = 0.9148771016999998
$ws.Range("S11").Value = 0.0000014257293286146
$ws.Range("T11").Value = 0.0000014257293286146

$ws.Range("I12").Value = 0.0002834193769566641
$ws.Range("J12").Value = 0.0002834193769566642
$ws.Range("M12").Value = 78.08909333333334
$ws.Range("N12").Value = 234.26728
$ws.Range("O12").Value = 0.4883264027331488
$ws.Range("P12").Value = 0.4883264027331488
$ws.Range("Q12").Value = 9.867858427555555
$ws.Range("R12").Value = 88.810725848
$ws.Range("S12").Value = 0.0001384011648141181
$ws.Range("T12").Value = 0.0001384011648141181

$ws.Range("I13").Value = 0.0002834193769566641
$ws.Range("J13").Value = 0.0002834193769566642
$ws.Range("M13").Value = 23.323433
$ws.Range("N13").Value = 69.970299
$ws.Range("O13").Value = 0.1458519704878668
$ws.Range("P13").Value = 0.1458519704878668
$ws.Range("Q13").Value = 2.947304483433333
$ws.Range("R13").Value = 26.5257403509
$ws.Range("S13").Value = 0.00004133727460357297
$ws.Range("T13").Value = 0.00004133727460357298

$ws.Range("G14").Value = 279.1165820000001
$ws.Range("H14").Value = 837.3497460000001
$ws.Range("I14").Value = 0.626011984453023
$ws.Range("J14").Value = 0.626011984453023
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 1.105124
$ws.Range("N14").Value = 3.315372
$ws.Range("O14").Value = 0.006910839970832482
$ws.Range("P14").Value = 0.006910839970832482
$ws.Range("Q14").Value = 308.4584335661681
$ws.Range("R14").Value = 2776.125902095512
$ws.Range("S14").Value = 0.004326268644378114
$ws.Range("T14").Value = 0.004326268644378114

$ws.Range("G15").Value = 279.1165820000001
$ws.Range("H15").Value = 837.3497460000001
$ws.Range("I15").Value = 0.626011984453023
$ws.Range("J15").Value = 0.626011984453023
$ws.Range("O15").Value = 0.0002777950170396876
$ws.Range("P15").Value = 0.0002777950170396876
$ws.Range("Q15").Value = 12.39910288332534
$ws.Range("R15").Value = 111.591925949928
$ws.Range("S15").Value = 0.0001739030098881762
$ws.Range("T15").Value = 0.0001739030098881762

$ws.Range("G16").Value = 279.1165820000001
$ws.Range("H16").Value = 837.3497460000001
$ws.Range("I16").Value = 0.626011984453023
$ws.Range("J16").Value = 0.626011984453023
$ws.Range("M16").Value = 56.54517366666666
$ws.Range("N16").Value = 169.635521
$ws.Range("O16").Value = 0.3536025335919447
$ws.Range("P16").Value = 0.3536025335919447
$ws.Range("Q16").Value = 15782.69560243641
$ws.Range("R16").Value = 142044.2604219277
$ws.Range("S16").Value = 0.22135942376151
$ws.Range("T16").Value = 0.22135942376151

$ws.Range("G17").Value = 279.1165820000001
$ws.Range("H17").Value = 837.3497460000001
$ws.Range("I17").Value = 0.626011984453023
$ws.Range("J17").Value = 0.626011984453023
$ws.Range("M17").Value = 0.8044289999999998
$ws.Range("N17").Value = 2.413287
$ws.Range("O17").Value = 0.005030458199167516
$ws.Range("P17").Value = 0.005030458199167516
$ws.Range("Q17").Value = 224.529472941678
$ws.Range("R17").Value = 2020.765256475102
$ws.Range("S17").Value = 0.003149127119968837
$ws.Range("T17").Value = 0.003149127119968837

$ws.Range("G18").Value = 279.1165820000001
$ws.Range("H18").Value = 837.3497460000001
$ws.Range("I18").Value = 0.626011984453023
$ws.Range("J18").Value = 0.626011984453023
$ws.Range("M18").Value = 78.08909333333334
$ws.Range("N18").Value = 234.26728
$ws.Range("O18").Value = 0.4883264027331488
$ws.Range("P18").Value = 0.4883264027331488
$ws.Range("Q18").Value = 21795.96082267899
$ws.Range("R18").Value = 196163.6474041109
$ws.Range("S18").Value = 0.3056981804357846
$ws.Range("T18").Value = 0.3056981804357846

$ws.Range("G19").Value = 279.1165820000001
$ws.Range("H19").Value = 837.3497460000001
$ws.Range("I19").Value = 0.626011984453023
$ws.Range("J19").Value = 0.626011984453023
$ws.Range("M19").Value = 23.323433
$ws.Range("N19").Value = 69.970299
$ws.Range("O19").Value = 0.1458519704878668
$ws.Range("P19").Value = 0.1458519704878668
$ws.Range("Q19").Value = 6509.956899466007
$ws.Range("R19").Value = 58589.61209519406
$ws.Range("S19").Value = 0.09130508148149326
$ws.Range("T19").Value = 0.09130508148149326

$ws.Range("G20").Value = 0.034934
$ws.Range("H20").Value = 0.104802
$ws.Range("I20").Value = 0.00007835114097550069
$ws.Range("J20").Value = 0.0000783511409755007
$ws.Range("K20").Value = 3.0
$ws.Range("L20").Value = 1.0
$ws.Range("M20").Value = 1.105124
$ws.Range("N20").Value = 3.315372
$ws.Range("O20").Value = 0.006910839970832482
$ws.Range("P20").Value = 0.006910839970832482
$ws.Range("Q20").Value = 0.038606401816
$ws.Range("R20").Value = 0.347457616344
$ws.Range("S20").Value = 0.000000541472196813821
$ws.Range("T20").Value = 0.000000541472196813821

$ws.Range("G21").Value = 0.034934
$ws.Range("H21").Value = 0.104802
$ws.Range("I21").Value = 0.00007835114097550069
$ws.Range("J21").Value = 0.0000783511409755007
$ws.Range("O21").Value = 0.0002777950170396876
$ws.Range("P21").Value = 0.0002777950170396876
$ws.Range("Q21").Value = 0.001551861437333333
$ws.Range("R21").Value = 0.013966752936
$ws.Range("S21").Value = 0.00000002176555654236818
$ws.Range("T21").Value = 0.00000002176555654236818

$ws.Range("G22").Value = 0.034934
$ws.Range("H22").Value = 0.104802
$ws.Range("I22").Value = 0.00007835114097550069
$ws.Range("J22").Value = 0.0000783511409755007
$ws.Range("M22").Value = 56.54517366666666
$ws.Range("N22").Value = 169.635521
$ws.Range("O22").Value = 0.3536025335919447
$ws.Range("P22").Value = 0.3536025335919447
$ws.Range("Q22").Value = 1.975349096871333
$ws.Range("R22").Value = 17.778141871842
$ws.Range("S22").Value = 0.00002770516195875667
$ws.Range("T22").Value = 0.00002770516195875668

$ws.Range("G23").Value = 0.034934
$ws.Range("H23").Value = 0.104802
$ws.Range("I23").Value = 0.00007835114097550069
$ws.Range("J23").Value = 0.0000783511409755007
$ws.Range("M23").Value = 0.8044289999999998
$ws.Range("N23").Value = 2.413287
$ws.Range("O23").Value = 0.005030458199167516
$ws.Range("P23").Value = 0.005030458199167516
$ws.Range("Q23").Value = 0.02810192268599999
$ws.Range("R23").Value = 0.252917304174
$ws.Range("S23").Value = 0.0000003941421395343373
$ws.Range("T23").Value = 0.0000003941421395343374

$ws.Range("G24").Value = 0.034934
$ws.Range("H24").Value = 0.104802
$ws.Range("I24").Value = 0.00007835114097550069
$ws.Range("J24").Value = 0.0000783511409755007
$ws.Range("M24").Value = 78.08909333333334
$ws.Range("N24").Value = 234.26728
$ws.Range("O24").Value = 0.4883264027331488
$ws.Range("P24").Value = 0.4883264027331488
$ws.Range("Q24").Value = 2.727964386506667
$ws.Range("R24").Value = 24.55167947856
$ws.Range("S24").Value = 0.00003826093082260407
$ws.Range("T24").Value = 0.00003826093082260407

$ws.Range("G25").Value = 0.034934
$ws.Range("H25").Value = 0.104802
$ws.Range("I25").Value = 0.00007835114097550069
$ws.Range("J25").Value = 0.0000783511409755007
$ws.Range("M25").Value = 23.323433
$ws.Range("N25").Value = 69.970299
$ws.Range("O25").Value = 0.1458519704878668
$ws.Range("P25").Value = 0.1458519704878668
$ws.Range("Q25").Value = 0.814780808422
$ws.Range("R25").Value = 7.333027275798
$ws.Range("S25").Value = 0.00001142766830124942
$ws.Range("T25").Value = 0.00001142766830124942

$ws.Range("G26").Value = 31.86718666666667
$ws.Range("H26").Value = 95.60156
$ws.Range("I26").Value = 0.07147278968948864
$ws.Range("J26").Value = 0.07147278968948864
$ws.Range("K26").Value = 3.0
$ws.Range("L26").Value = 1.0
$ws.Range("M26").Value = 1.105124
$ws.Range("N26").Value = 3.315372
$ws.Range("O26").Value = 0.006910839970832482
$ws.Range("P26").Value = 0.006910839970832482
$ws.Range("Q26").Value = 35.21719279781333
$ws.Range("R26").Value = 316.95473518032
$ws.Range("S26").Value = 0.0004939370118130219
$ws.Range("T26").Value = 0.0004939370118130219

$ws.Range("G27").Value = 31.86718666666667
$ws.Range("H27").Value = 95.60156
$ws.Range("I27").Value = 0.07147278968948864
$ws.Range("J27").Value = 0.07147278968948864
$ws.Range("O27").Value = 0.0002777950170396876
$ws.Range("P27").Value = 0.0002777950170396876
$ws.Range("Q27").Value = 1.415625410897778
$ws.Range("R27").Value = 12.74062869808
$ws.Range("S27").Value = 0.0000198547848296655
$ws.Range("T27").Value = 0.0000198547848296655

$ws.Range("G28").Value = 31.86718666666667
$ws.Range("H28").Value = 95.60156
$ws.Range("I28").Value = 0.07147278968948864
$ws.Range("J28").Value = 0.07147278968948864
$ws.Range("M28").Value = 56.54517366666666
$ws.Range("N28").Value = 169.635521
$ws.Range("O28").Value = 0.3536025335919447
$ws.Range("P28").Value = 0.3536025335919447
$ws.Range("Q28").Value = 1801.935604334751
$ws.Range("R28").Value = 16217.42043901276
$ws.Range("S28").Value = 0.0252729595170874
$ws.Range("T28").Value = 0.0252729595170874

$ws.Range("G29").Value = 31.86718666666667
$ws.Range("H29").Value = 95.60156
$ws.Range("I29").Value = 0.07147278968948864
$ws.Range("J29").Value = 0.07147278968948864
$ws.Range("M29").Value = 0.8044289999999998
$ws.Range("N29").Value = 2.413287
$ws.Range("O29").Value = 0.005030458199167516
$ws.Range("P29").Value = 0.005030458199167516
$ws.Range("Q29").Value = 25.63488910308
$ws.Range("R29").Value = 230.71400192772
$ws.Range("S29").Value = 0.0003595408809108636
$ws.Range("T29").Value = 0.0003595408809108636

$ws.Range("G30").Value = 31.86718666666667
$ws.Range("H30").Value = 95.60156
$ws.Range("I30").Value = 0.07147278968948864
$ws.Range("J30").Value = 0.07147278968948864
$ws.Range("M30").Value = 78.08909333333334
$ws.Range("N30").Value = 234.26728
$ws.Range("O30").Value = 0.4883264027331488
$ws.Range("P30").Value = 0.4883264027331488
$ws.Range("Q30").Value = 2488.479713884089
$ws.Range("R30").Value = 22396.3174249568
$ws.Range("S30").Value = 0.03490205028237087
$ws.Range("T30").Value = 0.03490205028237087

$ws.Range("G31").Value = 31.86718666666667
$ws.Range("H31").Value = 95.60156
$ws.Range("I31").Value = 0.07147278968948864
$ws.Range("J31").Value = 0.07147278968948864
$ws.Range("M31").Value = 23.323433
$ws.Range("N31").Value = 69.970299
$ws.Range("O31").Value = 0.1458519704878668
$ws.Range("P31").Value = 0.1458519704878668
$ws.Range("Q31").Value = 743.2521931184933
$ws.Range("R31").Value = 6689.26973806644
$ws.Range("S31").Value = 0.01042444721247681
$ws.Range("T31").Value = 0.01042444721247681

$ws.Range("G32").Value = 132.7532756666667
$ws.Range("H32").Value = 398.259827
$ws.Range("I32").Value = 0.2977434767481109
$ws.Range("J32").Value = 0.2977434767481109
$ws.Range("K32").Value = 3.0
$ws.Range("L32").Value = 1.0
$ws.Range("M32").Value = 1.105124
$ws.Range("N32").Value = 3.315372
$ws.Range("O32").Value = 0.006910839970832482
$ws.Range("P32").Value = 0.006910839970832482
$ws.Range("Q32").Value = 146.7088310178493
$ws.Range("R32").Value = 1320.379479160644
$ws.Range("S32").Value = 0.002057657520165477
$ws.Range("T32").Value = 0.002057657520165477

$ws.Range("G33").Value = 132.7532756666667
$ws.Range("H33").Value = 398.259827
$ws.Range("I33").Value = 0.2977434767481109
$ws.Range("J33").Value = 0.2977434767481109
$ws.Range("O33").Value = 0.0002777950170396876
$ws.Range("P33").Value = 0.0002777950170396876
$ws.Range("Q33").Value = 5.897254513848444
$ws.Range("R33").Value = 53.075290624636
$ws.Range("S33").Value = 0.0000827116541966973
$ws.Range("T33").Value = 0.0000827116541966973

$ws.Range("G34").Value = 132.7532756666667
$ws.Range("H34").Value = 398.259827
$ws.Range("I34").Value = 0.2977434767481109
$ws.Range("J34").Value = 0.2977434767481109
$ws.Range("M34").Value = 56.54517366666666
$ws.Range("N34").Value = 169.635521
$ws.Range("O34").Value = 0.3536025335919447
$ws.Range("P34").Value = 0.3536025335919447
$ws.Range("Q34").Value = 7506.55702739054
$ws.Range("R34").Value = 67559.01324651485
$ws.Range("S34").Value = 0.1052828477386063
$ws.Range("T34").Value = 0.1052828477386063

$ws.Range("G35").Value = 132.7532756666667
$ws.Range("H35").Value = 398.259827
$ws.Range("I35").Value = 0.2977434767481109
$ws.Range("J35").Value = 0.2977434767481109
$ws.Range("M35").Value = 0.8044289999999998
$ws.Range("N35").Value = 2.413287
$ws.Range("O35").Value = 0.005030458199167516
$ws.Range("P35").Value = 0.005030458199167516
$ws.Range("Q35").Value = 106.790584791261
$ws.Range("R35").Value = 961.1152631213488
$ws.Range("S35").Value = 0.001497786113856177
$ws.Range("T35").Value = 0.001497786113856177

$ws.Range("G36").Value = 132.7532756666667
$ws.Range("H36").Value = 398.259827
$ws.Range("I36").Value = 0.2977434767481109
$ws.Range("J36").Value = 0.2977434767481109
$ws.Range("M36").Value = 78.08909333333334
$ws.Range("N36").Value = 234.26728
$ws.Range("O36").Value = 0.4883264027331488
$ws.Range("P36").Value = 0.4883264027331488
$ws.Range("Q36").Value = 10366.58293384006
$ws.Range("R36").Value = 93299.24640456056
$ws.Range("S36").Value = 0.1453960009376659
$ws.Range("T36").Value = 0.1453960009376659

$ws.Range("G37").Value = 132.7532756666667
$ws.Range("H37").Value = 398.259827
$ws.Range("I37").Value = 0.2977434767481109
$ws.Range("J37").Value = 0.2977434767481109
$ws.Range("M37").Value = 23.323433
$ws.Range("N37").Value = 69.970299
$ws.Range("O37").Value = 0.1458519704878668
$ws.Range("P37").Value = 0.1458519704878668
$ws.Range("Q37").Value = 3096.26213054203
$ws.Range("R37").Value = 27866.35917487827
$ws.Range("S37").Value = 0.04342647278362034
$ws.Range("T37").Value = 0.04342647278362034

